$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching the style used by the
# existing header row (column H, "IP").
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I0 / IF columns, rows 2-72.
$data = @(
  @(7,7),
  @(8,8),
  @(8,8),
  @(7,7),
  @(7,8),
  @(8,8),
  @(8,8),
  @(8,8),
  @(6,6),
  @(7,7),
  @(9,9),
  @(5,5),
  @(7,8),
  @(6,7),
  @(5,6),
  @(6,6),
  @(8,8),
  @(7,7),
  @(7,7),
  @(6,7),
  @(6,6),
  @(7,7),
  @(7,7),
  @(8,8),
  @(7,7),
  @(7,7),
  @(8,8),
  @(8,8),
  @(7,7),
  @(8,8),
  @(7,7),
  @(8,8),
  @(8,8),
  @(7,7),
  @(6,6),
  @(8,9),
  @(7,7),
  @(6,8),
  @(8,9),
  @(7,8),
  @(8,8),
  @(8,8),
  @(8,8),
  @(10,10),
  @(9,9),
  @(9,9),
  @(9,9),
  @(7,7),
  @(8,8),
  @(6,6),
  @(7,8),
  @(9,9),
  @(8,8),
  @(7,7),
  @(8,9),
  @(7,7),
  @(7,8),
  @(9,9),
  @(8,8),
  @(8,8),
  @(8,8),
  @(4,5),
  @(7,7),
  @(5,5),
  @(7,7),
  @(7,7),
  @(7,7),
  @(4,4),
  @(7,7),
  @(2,2),
  @(4,4)
)

$row = 2
foreach ($pair in $data) {
  $ws.Cells.Item($row, 9).Value = $pair[0]
  $ws.Cells.Item($row, 10).Value = $pair[1]
  $row++
}
